# Commit: Update edited session - Cache Bust ID: 1760177226345o3f2no22s
#
# The scanner session log had a stale/duplicate log entry at row 12
# (Student ID 200852, logged 11:23:31) removed from the sheet, which
# shifts every subsequent row up by one (rows 13-115 -> 12-114), and the
# worksheet/tab was renamed from "General_Surgery" to "Session".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale log entry row; Excel shifts rows 13:115 up to 12:114.
$ws.Rows.Item(12).Delete()

# Rename the sheet/tab.
$ws.Name = "Session"
